# v4.13 estable | Flujo completo validado en entorno empresa: calculo ANS y descarga Drive
#
# Applies to REPOSITORIO_PEDIDOS_CERRADOS.xlsx:
#   1. Adds two new header columns Y1=AREA_OPERATIVA, Z1=SUBZONA (copying the
#      bold/bordered header style from X1) and blank Y/Z placeholder cells
#      for every existing data row.
#   2. Re-types the INSTALACION id in M2:M4 as zero-padded text (they were
#      being read back as floating point numbers, losing the leading zero).
#   3. Bumps DIAS_TRANSCURRIDOS (column T) by one full day across every
#      existing data row (ANS recalculated a day later).
#   4. Inserts a brand-new closed-order row (pedido 23398679) at row 10,
#      pushing the former row 10 (pedido 23499958) down to row 11 with its
#      own DIAS_TRANSCURRIDOS bumped by one day too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Forces a literal text cell (keeps leading zeros / long digit strings
    # from being re-interpreted as a float) and then clears the temporary
    # "@" number-format back to the sheet default so no stray style sticks
    # to the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Set-BlankTextCell($range) {
    # Materialises an empty (but present) inline-string cell, matching the
    # placeholder Y/Z columns added alongside AREA_OPERATIVA/SUBZONA.
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. New header columns Y1 / Z1 (copy X1's bold + bordered header style)
# ---------------------------------------------------------------------
$ws.Range("Y1").Value = "AREA_OPERATIVA"
$ws.Range("Z1").Value = "SUBZONA"
$ws.Range("X1").Copy()
$ws.Range("Y1:Z1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. M2:M4 -> zero-padded text instead of scientific-notation numbers
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("M2") "054926309110740000"
Set-TextValue $ws.Range("M3") "055224007200130201"
Set-TextValue $ws.Range("M4") "055224007200130301"

# ---------------------------------------------------------------------
# 3. DIAS_TRANSCURRIDOS (T2:T9) gains one day across the board
# ---------------------------------------------------------------------
$ws.Range("T2").Value = "18 días 09:36"
$ws.Range("T3").Value = "16 días 17:48"
$ws.Range("T4").Value = "16 días 17:50"
$ws.Range("T5").Value = "16 días 17:51"
$ws.Range("T6").Value = "14 días 13:27"
$ws.Range("T7").Value = "16 días 10:40"
$ws.Range("T8").Value = "17 días 15:54"
$ws.Range("T9").Value = "23 días 07:29"

# ---------------------------------------------------------------------
# 4. New blank AREA_OPERATIVA / SUBZONA cells for the existing rows 2-9
# ---------------------------------------------------------------------
foreach ($r in 2..9) {
    Set-BlankTextCell $ws.Range("Y" + $r)
    Set-BlankTextCell $ws.Range("Z" + $r)
}

# ---------------------------------------------------------------------
# 5. Insert the new closed order as row 10, pushing pedido 23499958 down
#    to row 11 (Excel shifts formatting + values automatically).
# ---------------------------------------------------------------------
$ws.Rows.Item(10).Insert()

Set-TextValue $ws.Range("A10") "23398679"
$ws.Range("B10").Value = "ENERES"
$ws.Range("C10").Value = "NUEVO"
$ws.Range("D10").Value = "ENEDOM"
$ws.Range("E10").Value = "27/03/2025 15:02"
$ws.Range("F10").Value = 45924.42361111111
$ws.Range("G10").Value = 43874896
$ws.Range("H10").Value = "LUZ EDITH MONTOYA"
$ws.Range("I10").Value = "SIN DATOS"
$ws.Range("J10").Value = 3107440550
$ws.Range("K10").Value = "RURAL_130023595050000002"
$ws.Range("L10").Value = "GUARNE"
Set-TextValue $ws.Range("M10") "130023595050000002"
$ws.Range("N10").Value = "Medellín"
$ws.Range("O10").Value = "ARTER"
$ws.Range("P10").Value = "Habilitación Viviendas Oriente"
$ws.Range("Q10").Value = "Rural"
$ws.Range("R10").Value = 8
$ws.Range("S10").Value = 45936.42361111111
$ws.Range("T10").Value = "27 días 10:10"
$ws.Range("U10").Value = "VENCIDO"
$ws.Range("V10").Value = "VENCIDO"
$ws.Range("W10").Value = "CERRADO"
$ws.Range("X10").Value = "Ejecutado en Campo"
Set-BlankTextCell $ws.Range("Y10")
Set-BlankTextCell $ws.Range("Z10")

# The old row 10 (pedido 23499958) is now row 11; keep its numbers/dates
# intact but apply the same +1 day ANS bump and the new blank placeholders.
$ws.Range("T11").Value = "49 días 09:04"
Set-BlankTextCell $ws.Range("Y11")
Set-BlankTextCell $ws.Range("Z11")
